$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.115.01"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "1.912.16"
$ws.Range("E3").Value = "  +2.31%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5063"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4078"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08352"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.47"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.106"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.79"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.95%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.398"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.24%  "
$ws.Range("D14").Value = "1.903.97"
$ws.Range("E14").Value = "  +2.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.232"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.22%  "
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.49"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.00%  "
$ws.Range("E18").Value = "  +2.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06504"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.63%  "
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.940"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.77%  "
$ws.Range("D23").Value = "30.133.97"
$ws.Range("E23").Value = "  +0.40%  "
$ws.Range("E24").Value = "  +2.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.190"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.54%  "
$ws.Range("D26").Value = "2.126.76"
$ws.Range("E26").Value = "  +2.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.79%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.07"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.281"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.75"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.98%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.145"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +10.38%  "
$ws.Range("E32").Value = "  +1.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.958"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.794"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.59%  "
$ws.Range("E35").Value = "  +1.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.358"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06394"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2152"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.91%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6538"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.196"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.640"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.84%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.38"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.212"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.39"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.83%  "
$ws.Range("E45").Value = "  +3.91%  "
$ws.Range("E46").Value = "  +10.46%  "
$ws.Range("E47").Value = "  -0.10%  "
$ws.Range("E48").Value = "  +0.91%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "122.21"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "79.03"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.138"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.91%  "
